$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data_B = New-Object 'object[,]' 24,3
$data_B[0,0] = 0.5207228070979113
$data_B[0,1] = 0.1469444347302797
$data_B[0,2] = 0.04310057380408949
$data_B[1,0] = 0.4760722566516904
$data_B[1,1] = 0.1467719285681817
$data_B[1,2] = 0.04087274509502947
$data_B[2,0] = 0.4487487078082211
$data_B[2,1] = 0.1466798740983712
$data_B[2,2] = 0.03949149556009246
$data_B[3,0] = 0.4376379893352009
$data_B[3,1] = 0.1466458751946291
$data_B[3,2] = 0.03892529582365967
$data_B[4,0] = 0.4357945240412846
$data_B[4,1] = 0.1466404427545882
$data_B[4,2] = 0.03883107854780832
$data_B[5,0] = 0.4485987672798899
$data_B[5,1] = 0.1466794013105925
$data_B[5,2] = 0.03948387303522338
$data_B[6,0] = 0.5053085920417573
$data_B[6,1] = 0.1468820932555701
$data_B[6,2] = 0.04233520651903433
$data_B[7,0] = 0.6172219749585395
$data_B[7,1] = 0.1473884334784188
$data_B[7,2] = 0.04781969422979415
$data_B[8,0] = 0.6998497859752604
$data_B[8,1] = 0.1478253487236572
$data_B[8,2] = 0.05178293112383159
$data_B[9,0] = 0.7375224846797437
$data_B[9,1] = 0.1480379069546203
$data_B[9,2] = 0.05357134085846837
$data_B[10,0] = 0.7517997740303599
$data_B[10,1] = 0.1481203561260358
$data_B[10,2] = 0.05424645764811231
$data_B[11,0] = 0.748724406753098
$data_B[11,1] = 0.1481025125646838
$data_B[11,2] = 0.05410115374912294
$data_B[12,0] = 0.7386968605707125
$data_B[12,1] = 0.1480446510106503
$data_B[12,2] = 0.05362692569678984
$data_B[13,0] = 0.7325561686285766
$data_B[13,1] = 0.1480094633299629
$data_B[13,2] = 0.05333617091611842
$data_B[14,0] = 0.6973894304878172
$data_B[14,1] = 0.1478117329956063
$data_B[14,2] = 0.0516657600962489
$data_B[15,0] = 0.6758369834069811
$data_B[15,1] = 0.1476939479521704
$data_B[15,2] = 0.05063728399396439
$data_B[16,0] = 0.6634486212123534
$data_B[16,1] = 0.1476275019351547
$data_B[16,2] = 0.05004437010266116
$data_B[17,0] = 0.6592555376444693
$data_B[17,1] = 0.1476052286963139
$data_B[17,2] = 0.04984338675895117
$data_B[18,0] = 0.6781304520174842
$data_B[18,1] = 0.1477063519373161
$data_B[18,2] = 0.05074690818761241
$data_B[19,0] = 0.7416418890229579
$data_B[19,1] = 0.1480615934247993
$data_B[19,2] = 0.05376627563030212
$data_B[20,0] = 0.7832167113577384
$data_B[20,1] = 0.1483051645904538
$data_B[20,2] = 0.05572726321869226
$data_B[21,0] = 0.7610216159645802
$data_B[21,1] = 0.1481741317677461
$data_B[21,2] = 0.05468178722914274
$data_B[22,0] = 0.6770935672446114
$data_B[22,1] = 0.1477007401384398
$data_B[22,2] = 0.05069735216216031
$data_B[23,0] = 0.5868735050169676
$data_B[23,1] = 0.1472399523984933
$data_B[23,2] = 0.04634755100340726
$ws.Range("B2:D25").Value = $data_B

$data_F = New-Object 'object[,]' 24,2
$data_F[0,0] = 0.8769724483368364
$data_F[0,1] = 0.002453848141947394
$data_F[1,0] = 0.8785232598701427
$data_F[1,1] = 0.002456241627494174
$data_F[2,0] = 0.8799975084742968
$data_F[2,1] = 0.002457790177935385
$data_F[3,0] = 0.8807295611691757
$data_F[3,1] = 0.002458441133364043
$data_F[4,0] = 0.8808590475767346
$data_F[4,1] = 0.00245855042812108
$data_F[5,0] = 0.8800068496121938
$data_F[5,1] = 0.002457798876104818
$data_F[6,0] = 0.8773988232358008
$data_F[6,1] = 0.002454657066581165
$data_F[7,0] = 0.8764271884825519
$data_F[7,1] = 0.002449119684124872
$data_F[8,0] = 0.8782410708210477
$data_F[8,1] = 0.002445427857540801
$data_F[9,0] = 0.8796157026952756
$data_F[9,1] = 0.002443829307024844
$data_F[10,0] = 0.8802152731424115
$data_F[10,1] = 0.002443235548013756
$data_F[11,0] = 0.8800826299927778
$data_F[11,1] = 0.002443362910528925
$data_F[12,0] = 0.8796634458325414
$data_F[12,1] = 0.002443780226483727
$data_F[13,0] = 0.8794169754249097
$data_F[13,1] = 0.002444037350475773
$data_F[14,0] = 0.8781622935055324
$data_F[14,1] = 0.002445533949598063
$data_F[15,0] = 0.8775333323565491
$data_F[15,1] = 0.002446472742172945
$data_F[16,0] = 0.8772232903068016
$data_F[16,1] = 0.002447020326832266
$data_F[17,0] = 0.8771271979466064
$data_F[17,1] = 0.002447207039032301
$data_F[18,0] = 0.8775949334546596
$data_F[18,1] = 0.002446372018292853
$data_F[19,0] = 0.8797844255204836
$data_F[19,1] = 0.002443657336694591
$data_F[20,0] = 0.8816760205874132
$data_F[20,1] = 0.002441950593348007
$data_F[21,0] = 0.8806242890337614
$data_F[21,1] = 0.00244285535951707
$data_F[22,0] = 0.8775669230049985
$data_F[22,1] = 0.00244641753109228
$data_F[23,0] = 0.8762463015492941
$data_F[23,1] = 0.002450551309323709
$ws.Range("F2:G25").Value = $data_F

$data_I = New-Object 'object[,]' 24,1
$data_I[0,0] = 0.8444350495859503
$data_I[1,0] = 0.8514222371491798
$data_I[2,0] = 0.8561503377242303
$data_I[3,0] = 0.8581871937662129
$data_I[4,0] = 0.8585320627433468
$data_I[5,0] = 0.8561773616286636
$data_I[6,0] = 0.8467533334220931
$data_I[7,0] = 0.8317481948379424
$data_I[8,0] = 0.8228435728691821
$data_I[9,0] = 0.8192531351800127
$data_I[10,0] = 0.8179597435858312
$data_I[11,0] = 0.8182353524759378
$data_I[12,0] = 0.8191453995382219
$data_I[13,0] = 0.8197114557021834
$data_I[14,0] = 0.8230874822915624
$data_I[15,0] = 0.8252764897093527
$data_I[16,0] = 0.8265788741980771
$data_I[17,0] = 0.8270272792786173
$data_I[18,0] = 0.8250389818633153
$data_I[19,0] = 0.8188762989291405
$data_I[20,0] = 0.8152346805192678
$data_I[21,0] = 0.8171429453125043
$data_I[22,0] = 0.8251462224321386
$data_I[23,0] = 0.8354353197679707
$ws.Range("I2:I25").Value = $data_I

$data_K = New-Object 'object[,]' 24,2
$data_K[0,0] = 0.3103589639591178
$data_K[0,1] = 0.3027314813407571
$data_K[1,0] = 0.2724810991288109
$data_K[1,1] = 0.2915323450551597
$data_K[2,0] = 0.2491835487336687
$data_K[2,1] = 0.2848281627347546
$data_K[3,0] = 0.2396800007927311
$data_K[3,1] = 0.2821395033813587
$data_K[4,0] = 0.2381013798646592
$data_K[4,1] = 0.2816956736681391
$data_K[5,0] = 0.2490554187315297
$data_K[5,1] = 0.2847917269269544
$data_K[6,0] = 0.2973073760006173
$data_K[6,1] = 0.2988343280366763
$data_K[7,0] = 0.3915893237252419
$data_K[7,1] = 0.3277365301812267
$data_K[8,0] = 0.4606315092820239
$data_K[8,1] = 0.3498040966372855
$data_K[9,0] = 0.491987649049122
$data_K[9,1] = 0.3600245084145115
$data_K[10,0] = 0.5038535365653729
$data_K[10,1] = 0.3639208223555386
$data_K[11,0] = 0.5012983707581782
$data_K[11,1] = 0.3630805233028269
$data_K[12,0] = 0.4929640265421824
$data_K[12,1] = 0.3603445385239468
$data_K[13,0] = 0.4878579402716809
$data_K[13,1] = 0.3586720612397727
$data_K[14,0] = 0.4585812244832823
$data_K[14,1] = 0.3491398201537663
$data_K[15,0] = 0.4406072926653621
$data_K[15,1] = 0.3433386026024152
$data_K[16,0] = 0.4302643509061284
$data_K[16,1] = 0.34001900864709
$data_K[17,0] = 0.4267616003742774
$data_K[17,1] = 0.3388979925086062
$data_K[18,0] = 0.442521150701964
$data_K[18,1] = 0.3439543815071886
$data_K[19,0] = 0.4954122473872076
$data_K[19,1] = 0.3611474567328656
$data_K[20,0] = 0.5299327187330221
$data_K[20,1] = 0.3725360228545753
$data_K[21,0] = 0.5115129831730485
$data_K[21,1] = 0.3664438560922179
$data_K[22,0] = 0.4416559249968941
$data_K[22,1] = 0.3436759392587305
$data_K[23,0] = 0.3661219813189689
$data_K[23,1] = 0.3197714952218007
$ws.Range("K2:L25").Value = $data_K

$data_O = New-Object 'object[,]' 24,1
$data_O[0,0] = 3.106186664742836
$data_O[1,0] = 3.124885900227142
$data_O[2,0] = 3.138093725426785
$data_O[3,0] = 3.143910097293542
$data_O[4,0] = 3.144902117184103
$data_O[5,0] = 3.13817040956053
$data_O[6,0] = 3.112275812109075
$data_O[7,0] = 3.075198868358882
$data_O[8,0] = 3.056319809934251
$data_O[9,0] = 3.049548765482399
$data_O[10,0] = 3.047246170471198
$data_O[11,0] = 3.047730444948797
$data_O[12,0] = 3.049354088848588
$data_O[13,0] = 3.050382670908192
$data_O[14,0] = 3.056798891937689
$data_O[15,0] = 3.061200555297859
$data_O[16,0] = 3.063903303970221
$data_O[17,0] = 3.064847777888559
$data_O[18,0] = 3.060714288825437
$data_O[19,0] = 3.048870088336258
$data_O[20,0] = 3.042653233006632
$data_O[21,0] = 3.045831788906384
$data_O[22,0] = 3.060933593594854
$data_O[23,0] = 3.08376139890197
$ws.Range("O2:O25").Value = $data_O
